# Commit: "Se agrega nuevo componente"
# A new DAO component (AutorizacionTipoDAO) was inserted earlier in the
# "Daos" progress list, which pushes CategoriaAdquisicionDAO and
# ColaboradorDAO down by one row, and the ComponenteTipoDAO component is
# now marked complete (100%).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daos")

$ws.Range("B10").Value = "AutorizacionTipoDAO"
$ws.Range("C10").Value = 1

$ws.Range("B11").Value = "CategoriaAdquisicionDAO"
$ws.Range("C11").Value = 0

$ws.Range("B12").Value = "ColaboradorDAO"
$ws.Range("C12").Value = 0

$ws.Range("C17").Value = 1

# Reflect the resulting selection/view state.
$ws.Activate()
$ws.Range("C18").Select()
